$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.626.49'
$ws.Range('E2').Value = '  +0.84%  '
$ws.Range('D3').Value = '2.240.27'
$ws.Range('E3').Value = '  -0.01%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').Value = "'305.40"
$ws.Range('E5').Value = '  -0.51%  '
$ws.Range('D6').Value = "'93.29"
$ws.Range('E6').Value = '  -1.28%  '
$ws.Range('D7').Value = "'0.567"
$ws.Range('E7').Value = '  -0.55%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').Value = "'0.511"
$ws.Range('E9').Value = '  -2.53%  '
$ws.Range('D10').Value = "'34.53"
$ws.Range('E10').Value = '  -0.06%  '
$ws.Range('D11').Value = "'0.0798"
$ws.Range('E11').Value = '  -1.60%  '
$ws.Range('D12').Value = "'7.12"
$ws.Range('E12').Value = '  -0.51%  '
$ws.Range('D13').Value = "'0.104"
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').Value = '2.586.50'
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Value = '2.339.46'
$ws.Range('E15').Value = '  +0.17%  '
$ws.Range('D16').Value = "'0.825"
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').Value = "'13.46"
$ws.Range('E17').Value = '  -0.23%  '
$ws.Range('D18').Value = '44.405.02'
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').Value = '0.0₃0930'
$ws.Range('E19').Value = '  -3.50%  '
$ws.Range('D20').Value = "'6.14"
$ws.Range('E20').Value = '  -3.87%  '
$ws.Range('D21').Value = "'11.64"
$ws.Range('E21').Value = '  -4.29%  '
$ws.Range('D22').Value = "'65.18"
$ws.Range('E22').Value = '  -0.49%  '
$ws.Range('D23').Value = "'236.77"
$ws.Range('E23').Value = '  -0.23%  '
$ws.Range('D24').Value = "'2.93"
$ws.Range('E24').Value = '  -0.59%  '
$ws.Range('D25').Value = "'1.95"
$ws.Range('E25').Value = '  -1.90%  '
$ws.Range('E26').Value = '  -0.10%  '
$ws.Range('E27').Value = '  +3.59%  '
$ws.Range('D28').Value = "'9.71"
$ws.Range('E28').Value = '  -1.96%  '
$ws.Range('D29').Value = "'36.75"
$ws.Range('E29').Value = '  -4.92%  '
$ws.Range('D30').Value = "'19.86"
$ws.Range('E30').Value = '  -0.87%  '
$ws.Range('D31').Value = "'5.78"
$ws.Range('E31').Value = '  -1.01%  '
$ws.Range('D32').Value = "'148.52"
$ws.Range('E32').Value = '  -3.24%  '
$ws.Range('E33').Value = '  +0.30%  '
$ws.Range('D34').Value = "'0.0774"
$ws.Range('E34').Value = '  -2.33%  '
$ws.Range('D35').Value = "'3.14"
$ws.Range('E35').Value = '  +0.07%  '
$ws.Range('E36').Value = '  +1.40%  '
$ws.Range('D37').Value = "'0.117"
$ws.Range('E37').Value = '  -2.26%  '
$ws.Range('D38').Value = "'1.84"
$ws.Range('E38').Value = '  +4.95%  '
$ws.Range('D39').Value = "'14.94"
$ws.Range('E39').Value = '  +4.84%  '
$ws.Range('D40').Value = "'3.32"
$ws.Range('E40').Value = '  -5.05%  '
$ws.Range('D41').Value = "'3.75"
$ws.Range('E41').Value = '  -1.17%  '
$ws.Range('D42').Value = "'0.0296"
$ws.Range('E42').Value = '  -0.83%  '
$ws.Range('E43').Value = '  +0.00%  '
$ws.Range('D44').Value = '1.806.45'
$ws.Range('E44').Value = '  +3.46%  '
$ws.Range('D45').Value = "'1.76"
$ws.Range('E45').Value = '  +10.73%  '
$ws.Range('D46').Value = "'81.07"
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').Value = "'0.185"
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('D48').Value = "'97.35"
$ws.Range('E48').Value = '  -2.17%  '
$ws.Range('D49').Value = "'4.81"
$ws.Range('E49').Value = '  -2.24%  '
$ws.Range('D50').Value = "'68.45"
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('D51').Value = "'53.48"
$ws.Range('E51').Value = '  -2.23%  '
